$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: Insert a new column before P, shifting old P (Hours) -> Q, old Q (Time) -> R ---
$ws.Columns("P").Insert()

# --- Step 2: Insert a new row before the existing total row (32), pushing it down to row 33 ---
$ws.Rows("32").Insert()

# --- Step 3: Jot down the "Reworking flocking" note (also backfilled onto row 20) ---
$ws.Range("E21").Value = "Reworking flocking"
$ws.Range("E20").Value = "Reworking flocking"

# --- Step 4: Fill in the new "Learning goal 2" mini-table (rows 3-7, columns N/O/Q) ---
$ws.Range("N3").NumberFormat = "d-mmm"
$ws.Range("N3").Value = 45091
$ws.Range("O3").Value = "18.00 - 22.00"
$ws.Range("Q3").Value = 4

$ws.Range("N4").NumberFormat = "d-mmm"
$ws.Range("N4").Value = 45098
$ws.Range("O4").Value = "18.00 - 22.00"
$ws.Range("Q4").Value = 4

# --- Step 5: Note the upcoming deadline ---
$ws.Range("J37").Value = "July 12th deadline"

$ws.Range("N5").NumberFormat = "d-mmm"
$ws.Range("N5").Value = 45102
$ws.Range("O5").Value = "11.00 - 15.00"
$ws.Range("Q5").Value = 4

$ws.Range("N6").NumberFormat = "d-mmm"
$ws.Range("N6").Value = 45105
$ws.Range("O6").Value = "18.00 - 22.00"
$ws.Range("Q6").Value = 4

$ws.Range("N7").NumberFormat = "d-mmm"
$ws.Range("N7").Value = 45112
$ws.Range("O7").Value = "18.00 - 22.00"
$ws.Range("Q7").Value = 4

# --- Step 6: Combined grand-total label ---
$ws.Range("N37").Value = "Total hours both goals = "

# --- Step 7: "Learning goal 2" subtotal row (32) ---
$ws.Range("O32").Value = "Total hours="
$ws.Range("Q32").Formula = "=SUM(Q3:Q31)"

# --- Step 8: Finish row 21 (primary table) ---
$ws.Range("B21").Value = "21.00 - 22.15"
$ws.Range("A21").NumberFormat = "d-mmm"
$ws.Range("A21").Value = 45069
$ws.Range("D21").Value = 1.15

# --- Step 9: Row 22 (primary table) ---
$ws.Range("E22").Value = "Reimplemented flocking. Improved neighbour finding indicators."
$ws.Range("A22").NumberFormat = "d-mmm"
$ws.Range("A22").Value = 45070
$ws.Range("B22").Value = "18.00 - 22.00"
$ws.Range("D22").Value = 4

# --- Step 10: Rows 23-26 (primary table) ---
$ws.Range("A23").NumberFormat = "d-mmm"
$ws.Range("A23").Value = 45074
$ws.Range("B23").Value = "11.00 - 15.00"
$ws.Range("D23").Value = 4

$ws.Range("A24").NumberFormat = "d-mmm"
$ws.Range("A24").Value = 45077
$ws.Range("B24").Value = "18.00 - 22.00"
$ws.Range("D24").Value = 4

$ws.Range("A25").NumberFormat = "d-mmm"
$ws.Range("A25").Value = 45081
$ws.Range("B25").Value = "11.00 - 15.00"
$ws.Range("D25").Value = 4

$ws.Range("A26").NumberFormat = "d-mmm"
$ws.Range("A26").Value = 45084
$ws.Range("B26").Value = "18.00 - 22.00"
$ws.Range("D26").Value = 4

# --- Step 11: Update the "Learning goal 1" total row (now row 33) to include the new rows ---
$ws.Range("D33").Formula = "=SUM(D3:D32)"

# --- Step 12: Combined grand total formula ---
$ws.Range("Q37").Formula = "=SUM(D33,Q32)"

# --- Step 13: Update the selected cell to match the saved workbook state ---
$ws.Range("G29").Select()
